$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values to add on row 2 (stored as text, matching the source inlineStr cells)
$values = @("07/07/2023", "1000.00", "1000.00", "1000.00", "1000.00", "0.00", "100.00")

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(2, $col)
    # Prefix with an apostrophe so Excel stores the value as literal text
    # (avoids auto-conversion of dates/numbers), then reset the cell style
    # back to Normal so no extra style index / number format is introduced.
    $cell.Value = "'" + $values[$i]
    $cell.Style = "Normal"
}
